$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 129
$ws.Range("H129").Value = 1036.5306
$ws.Range("I129").Value = 570
$ws.Range("J129").Value = 1056.3829
$ws.Range("K129").Value = 1710
$ws.Range("L129").Value = 3169.1487
$ws.Range("M129").Value = 3290
$ws.Range("N129").Value = -13169.1487
# Row 137
$ws.Range("H137").Value = 2036.1471
$ws.Range("I137").Value = 2151.2917
$ws.Range("J137").Value = 1759.8
$ws.Range("K137").Value = 6453.875100000001
$ws.Range("L137").Value = 5279.4
$ws.Range("M137").Value = -3903.875100000001
$ws.Range("N137").Value = -10379.4

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 1196.0667
$ws.Range("I2").Value = 1295.0385
$ws.Range("J2").Value = 552.75
$ws.Range("K2").Value = 1295.0385
$ws.Range("L2").Value = 552.75
$ws.Range("M2").Value = -1182.0385
$ws.Range("N2").Value = -778.75
# Row 32
$ws.Range("H32").Value = 28674.35
$ws.Range("I32").Value = 30651.666
$ws.Range("J32").Value = 10878.5
$ws.Range("K32").Value = 30651.666
$ws.Range("L32").Value = 10878.5
$ws.Range("M32").Value = -30364.666
$ws.Range("N32").Value = -11452.5
# Row 61
$ws.Range("H61").Value = 5136.51
$ws.Range("I61").Value = 3737.1538
$ws.Range("J61").Value = 9684.416999999999
$ws.Range("K61").Value = 3737.1538
$ws.Range("L61").Value = 9684.416999999999
$ws.Range("M61").Value = -3525.1538
$ws.Range("N61").Value = -10108.417
# Row 74
$ws.Range("H74").Value = 3526.2979
$ws.Range("I74").Value = 1407.6279
$ws.Range("J74").Value = 26302
$ws.Range("K74").Value = 1407.6279
$ws.Range("L74").Value = 26302
$ws.Range("M74").Value = -533.6279
$ws.Range("N74").Value = -28050
# Row 77
$ws.Range("H77").Value = 3526.2979
$ws.Range("I77").Value = 1407.6279
$ws.Range("J77").Value = 26302
$ws.Range("K77").Value = 7038.139499999999
$ws.Range("L77").Value = 131510
$ws.Range("M77").Value = -2670.139499999999
$ws.Range("N77").Value = -140246
# Row 116
$ws.Range("H116").Value = 1196.0667
$ws.Range("I116").Value = 1295.0385
$ws.Range("J116").Value = 552.75
$ws.Range("K116").Value = 1295.0385
$ws.Range("L116").Value = 552.75
$ws.Range("M116").Value = 998.9614999999999
$ws.Range("N116").Value = -5140.75
# Row 122
$ws.Range("H122").Value = 4809708.5
$ws.Range("I122").Value = 2128.762
$ws.Range("J122").Value = 25001542
$ws.Range("K122").Value = 6386.286
$ws.Range("L122").Value = 75004626
$ws.Range("M122").Value = -3936.286
$ws.Range("N122").Value = -75009526
# Row 136
$ws.Range("H136").Value = 5136.51
$ws.Range("I136").Value = 3737.1538
$ws.Range("J136").Value = 9684.416999999999
$ws.Range("K136").Value = 11211.4614
$ws.Range("L136").Value = 29053.251
$ws.Range("M136").Value = -8661.4614
$ws.Range("N136").Value = -34153.251

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 1196.0667
$ws.Range("I3").Value = 1295.0385
$ws.Range("J3").Value = 552.75
$ws.Range("K3").Value = 1295.0385
$ws.Range("L3").Value = 552.75
$ws.Range("M3").Value = -1181.0385
$ws.Range("N3").Value = -780.75
# Row 31
$ws.Range("H31").Value = 17750
$ws.Range("I31").Value = 500
$ws.Range("K31").Value = 500
$ws.Range("M31").Value = -248
# Row 62
$ws.Range("H62").Value = 40181
$ws.Range("J62").Value = 40181
$ws.Range("L62").Value = 40181
$ws.Range("N62").Value = -41553
# Row 65
$ws.Range("H65").Value = 40181
$ws.Range("J65").Value = 40181
$ws.Range("L65").Value = 120543
$ws.Range("N65").Value = -127407
# Row 86
$ws.Range("H86").Value = 6412051.5
$ws.Range("I86").Value = 7093976
$ws.Range("J86").Value = 1961.4
$ws.Range("K86").Value = 7093976
$ws.Range("L86").Value = 1961.4
$ws.Range("M86").Value = -7092853
$ws.Range("N86").Value = -4207.4
# Row 89
$ws.Range("H89").Value = 6412051.5
$ws.Range("I89").Value = 7093976
$ws.Range("J89").Value = 1961.4
$ws.Range("K89").Value = 35469880
$ws.Range("L89").Value = 9807
$ws.Range("M89").Value = -35464264
$ws.Range("N89").Value = -21039
# Row 134
$ws.Range("H134").Value = 1654.9
$ws.Range("I134").Value = 1505.8334
$ws.Range("J134").Value = 1878.5
$ws.Range("K134").Value = 4517.5002
$ws.Range("L134").Value = 5635.5
$ws.Range("M134").Value = -1982.5002
$ws.Range("N134").Value = -10705.5

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 4
$ws.Range("H4").Value = 299.5
$ws.Range("J4").Value = 299.5
$ws.Range("L4").Value = 299.5
$ws.Range("N4").Value = -523.5
# Row 31
$ws.Range("H31").Value = 2528.2632
$ws.Range("I31").Value = 1557.6111
$ws.Range("J31").Value = 20000
$ws.Range("K31").Value = 1557.6111
$ws.Range("L31").Value = 20000
$ws.Range("M31").Value = -1262.6111
$ws.Range("N31").Value = -20590
# Row 34
$ws.Range("H34").Value = 2528.2632
$ws.Range("I34").Value = 1557.6111
$ws.Range("J34").Value = 20000
$ws.Range("K34").Value = 1557.6111
$ws.Range("L34").Value = 20000
$ws.Range("M34").Value = -1355.6111
$ws.Range("N34").Value = -20404
# Row 58
$ws.Range("H58").Value = 2599300.5
$ws.Range("I58").Value = 3637780.5
$ws.Range("K58").Value = 3637780.5
$ws.Range("M58").Value = -3637577.5
# Row 132
$ws.Range("H132").Value = 2013.76
$ws.Range("I132").Value = 1832.7441
$ws.Range("J132").Value = 3125.7144
$ws.Range("K132").Value = 5498.2323
$ws.Range("L132").Value = 9377.143199999999
$ws.Range("M132").Value = -2968.2323
$ws.Range("N132").Value = -14437.1432
# Row 134
$ws.Range("H134").Value = 2541.8333
$ws.Range("I134").Value = 1965.8
$ws.Range("J134").Value = 5422
$ws.Range("K134").Value = 5897.4
$ws.Range("L134").Value = 16266
$ws.Range("M134").Value = -3362.4
$ws.Range("N134").Value = -21336
# Row 136
$ws.Range("H136").Value = 2599300.5
$ws.Range("I136").Value = 3637780.5
$ws.Range("K136").Value = 10913341.5
$ws.Range("M136").Value = -10910791.5

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 4
$ws.Range("H4").Value = 1500
$ws.Range("I4").Value = 1000
$ws.Range("J4").Value = 2000
$ws.Range("K4").Value = 3000
$ws.Range("L4").Value = 6000
$ws.Range("M4").Value = -2888
$ws.Range("N4").Value = -6224
# Row 70
$ws.Range("H70").Value = 3401
$ws.Range("I70").Value = 1306
$ws.Range("K70").Value = 3918
$ws.Range("M70").Value = -3603
# Row 73
$ws.Range("H73").Value = 3401
$ws.Range("I73").Value = 1306
$ws.Range("K73").Value = 3918
$ws.Range("M73").Value = -2826
# Row 75
$ws.Range("H75").Value = 3672.3635
$ws.Range("I75").Value = 1750
$ws.Range("J75").Value = 4099.5557
$ws.Range("K75").Value = 5250
$ws.Range("L75").Value = 12298.6671
$ws.Range("M75").Value = -4252
$ws.Range("N75").Value = -14294.6671
# Row 78
$ws.Range("H78").Value = 3672.3635
$ws.Range("I78").Value = 1750
$ws.Range("J78").Value = 4099.5557
$ws.Range("K78").Value = 15750
$ws.Range("L78").Value = 36896.0013
$ws.Range("M78").Value = -10758
$ws.Range("N78").Value = -46880.0013
# Row 98
$ws.Range("H98").Value = 469.66666
$ws.Range("I98").Value = 392.94116
$ws.Range("K98").Value = 1178.82348
$ws.Range("M98").Value = 319.17652

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 5
$ws.Range("H5").Value = 19666.666
$ws.Range("J5").Value = 19666.666
$ws.Range("L5").Value = 19666.666
$ws.Range("N5").Value = -19890.666
# Row 102
$ws.Range("H102").Value = 2805.0605
$ws.Range("I102").Value = 2344.1924
$ws.Range("J102").Value = 4516.857
$ws.Range("K102").Value = 2344.1924
$ws.Range("L102").Value = 4516.857
$ws.Range("M102").Value = -722.1923999999999
$ws.Range("N102").Value = -7760.857
# Row 113
$ws.Range("H113").Value = 2495
$ws.Range("I113").Value = 2560.7334
$ws.Range("J113").Value = 2166.3333
$ws.Range("K113").Value = 2560.7334
$ws.Range("L113").Value = 2166.3333
$ws.Range("M113").Value = -390.7334000000001
$ws.Range("N113").Value = -6506.3333
# Row 126
$ws.Range("H126").Value = 2447.3635
$ws.Range("I126").Value = 1884
$ws.Range("J126").Value = 3261.111
$ws.Range("K126").Value = 5652
$ws.Range("L126").Value = 9783.332999999999
$ws.Range("M126").Value = -3182
$ws.Range("N126").Value = -14723.333

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 1300.6666
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 1300.6666
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 1300.6666
$ws.Range("M22").ClearContents()
$ws.Range("N22").Value = -1890.6666
# Row 27
$ws.Range("H27").Value = 1300.6666
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 1300.6666
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 1300.6666
$ws.Range("M27").ClearContents()
$ws.Range("N27").Value = -1514.6666
# Row 132
$ws.Range("H132").Value = 2915.9119
$ws.Range("I132").Value = 2236.3462
$ws.Range("J132").Value = 5124.5
$ws.Range("K132").Value = 6709.0386
$ws.Range("L132").Value = 15373.5
$ws.Range("M132").Value = -4179.0386
$ws.Range("N132").Value = -20433.5
# Row 136
$ws.Range("H136").Value = 4598.341
$ws.Range("I136").Value = 2724.2917
$ws.Range("J136").Value = 6847.2
$ws.Range("K136").Value = 8172.875100000001
$ws.Range("L136").Value = 20541.6
$ws.Range("M136").Value = -5622.875100000001
$ws.Range("N136").Value = -25641.6

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 122
$ws.Range("H122").Value = 793.0833
$ws.Range("I122").Value = 821.2
$ws.Range("J122").Value = 652.5
$ws.Range("K122").Value = 2463.6
$ws.Range("L122").Value = 1957.5
$ws.Range("M122").Value = -13.60000000000036
$ws.Range("N122").Value = -6857.5
# Row 126
$ws.Range("H126").Value = 1375.2916
$ws.Range("I126").Value = 1410.35
$ws.Range("J126").Value = 1200
$ws.Range("K126").Value = 4231.049999999999
$ws.Range("L126").Value = 3600
$ws.Range("M126").Value = -1761.049999999999
$ws.Range("N126").Value = -8540
# Row 132
$ws.Range("H132").Value = 3102.182
$ws.Range("I132").Value = 2827.4
$ws.Range("J132").Value = 5850
$ws.Range("K132").Value = 8482.200000000001
$ws.Range("L132").Value = 17550
$ws.Range("M132").Value = -5952.200000000001
$ws.Range("N132").Value = -22610
# Row 136
$ws.Range("H136").Value = 4876.409
$ws.Range("I136").Value = 1657.48
$ws.Range("J136").Value = 9111.842000000001
$ws.Range("K136").Value = 4972.440000000001
$ws.Range("L136").Value = 27335.526
$ws.Range("M136").Value = -2422.440000000001
$ws.Range("N136").Value = -32435.526
